$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Helper: write $value into $range as TEXT (never auto-coerced to a
# number/percentage/date by Excel's smart entry) while leaving the
# cell's style at the default (no explicit number-format xf record
# left behind). We do this by temporarily forcing a text number
# format, assigning the value, then pasting the (default) formats
# back in from an untouched scratch cell on the same sheet.
# -----------------------------------------------------------------
function Set-TextValue {
    param($range, $value, $donor)

    $range.NumberFormat = "@"
    $range.Value = $value
    $donor.Copy()
    $range.PasteSpecial(-4122)  # xlPasteFormats
}

# --- Remove the stray empty inline-string value in "ODI Batting"!B2 ---
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$odiBatting.Range("B2").ClearContents()

# --- Add "ODI Batting Extra" sheet after "ODI Bowling" ---
$odiBowling = $wb.Worksheets.Item("ODI Bowling")
$battingExtra = $wb.Worksheets.Add($null, $odiBowling)
$battingExtra.Name = "ODI Batting Extra"
$beDonor = $battingExtra.Cells.Item(500, 500)

$battingExtra.Range("A1").Value = "MATCH_CODE"
$battingExtra.Range("B1").Value = "BATTING_POSITION"
$battingExtra.Range("C1").Value = "NUM_4"
$battingExtra.Range("D1").Value = "NUM_6"
$battingExtra.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$battingExtra.Range("F1").Value = "MAN_OF_MATCH"
$odiBowling.Range("A1").Copy()
$battingExtra.Range("A1:F1").PasteSpecial(-4122)  # xlPasteFormats (bold header style)

Set-TextValue $battingExtra.Range("A2") "4259" $beDonor
Set-TextValue $battingExtra.Range("F2") "NO" $beDonor

Set-TextValue $battingExtra.Range("A3") "4348" $beDonor
$battingExtra.Range("B3").Value = 10
Set-TextValue $battingExtra.Range("C3") "2" $beDonor
Set-TextValue $battingExtra.Range("D3") "2" $beDonor
Set-TextValue $battingExtra.Range("E3") "8.68%" $beDonor
Set-TextValue $battingExtra.Range("F3") "NO" $beDonor

# --- Add "ODI Bowling Extra" sheet after "ODI Batting Extra" ---
$bowlingExtra = $wb.Worksheets.Add($null, $battingExtra)
$bowlingExtra.Name = "ODI Bowling Extra"
$weDonor = $bowlingExtra.Cells.Item(500, 500)

$bowlingExtra.Range("A1").Value = "MATCH_CODE"
$bowlingExtra.Range("B1").Value = "MAIDEN_OVERS"
$bowlingExtra.Range("C1").Value = "PERCENT_WICKETS_OF_ALL"
$odiBowling.Range("A1").Copy()
$bowlingExtra.Range("A1:C1").PasteSpecial(-4122)  # xlPasteFormats (bold header style)

Set-TextValue $bowlingExtra.Range("A2") "4348" $weDonor
Set-TextValue $bowlingExtra.Range("B2") "0" $weDonor
Set-TextValue $bowlingExtra.Range("C2") "10.00%" $weDonor
